$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column AB (which currently holds "eFG%"),
# shifting eFG%..MATCHES (and their data) one column to the right.
$ws.Range("AB:AB").Insert()

# New header for the inserted column
$ws.Range("AB1").Value = "PPSA"

# New PPSA values for rows 2-12
$ws.Range("AB2").Value = 2.301533219761499
$ws.Range("AB3").Value = 2.286376274328082
$ws.Range("AB4").Value = 2.305084745762712
$ws.Range("AB5").Value = 2.240909090909091
$ws.Range("AB6").Value = 2.215189873417721
$ws.Range("AB7").Value = 2.288775510204081
$ws.Range("AB8").Value = 2.354948805460751
$ws.Range("AB9").Value = 2.352660841938046
$ws.Range("AB10").Value = 2.269799825935596
$ws.Range("AB11").Value = 2.31390134529148
$ws.Range("AB12").Value = 2.343525179856115

# Rename the TO ratio stat columns (now shifted to AH/AI)
$ws.Range("AH1").Value = "TOR%"
$ws.Range("AI1").Value = "OppTOR%"
